$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 606, shifting existing rows 606-659 down to 607-660
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row 606 with the new record
$ws.Cells.Item(606, 1).Value = 5
$ws.Cells.Item(606, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(606, 3).Value = "Maule"
$ws.Cells.Item(606, 4).Value = 45223
$ws.Cells.Item(606, 4).Style = $ws.Cells.Item(607, 4).Style
$ws.Cells.Item(606, 4).NumberFormat = $ws.Cells.Item(607, 4).NumberFormat
$ws.Cells.Item(606, 5).Value = 7
$ws.Cells.Item(606, 6).Value = 100112032
$ws.Cells.Item(606, 7).Value = "Zapallo italiano"
$ws.Cells.Item(606, 8).Value = "Sin especificar"
$ws.Cells.Item(606, 9).Value = "Primera"
$ws.Cells.Item(606, 10).Value = 200
$ws.Cells.Item(606, 11).Value = 15000
$ws.Cells.Item(606, 12).Value = 15000
$ws.Cells.Item(606, 13).Value = 15000
$ws.Cells.Item(606, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(606, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(606, 16).Value = 300
$ws.Cells.Item(606, 17).Value = 50
$ws.Cells.Item(606, 18).Value = "Hortaliza"
